$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: A2:C2 become numeric 1's (previously blank), D2 gets the new note text.
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "← Não mexer nessa linha"

# New column D formatting: bold font + wrap text so the note is readable.
$ws.Range("D2").Font.Bold = $true
$ws.Range("D2").WrapText = $true

# Row height grows to fit the wrapped text; column D gets a wider custom width
# (21.6 is the input that rounds to the target stored width of ~22.43/22.5).
$ws.Rows.Item(2).RowHeight = 23.85
$ws.Columns.Item(4).ColumnWidth = 21.6

# Selection moves from E11 to D12.
$ws.Range("D12").Select() | Out-Null

Write-Output "done"
